$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value2 = 2.0
$ws.Range("G2").Value2 = 22.556342
$ws.Range("H2").Value2 = 45.112684
$ws.Range("I2").Value2 = 0.0667629019027735
$ws.Range("J2").Value2 = 0.04665728030990886
$ws.Range("K2").Value2 = 2.0
$ws.Range("M2").Value2 = 29.4426765
$ws.Range("N2").Value2 = 58.88535299999999
$ws.Range("O2").Value2 = 0.2070274275189754
$ws.Range("P2").Value2 = 0.1657617245498156
$ws.Range("Q2").Value2 = 664.1190805293629
$ws.Range("R2").Value2 = 2656.476322117452
$ws.Range("S2").Value2 = 0.01382175183463291
$ws.Range("T2").Value2 = 0.007733991246974647

$ws.Range("E3").Value2 = 2.0
$ws.Range("G3").Value2 = 22.556342
$ws.Range("H3").Value2 = 45.112684
$ws.Range("I3").Value2 = 0.0667629019027735
$ws.Range("J3").Value2 = 0.04665728030990886
$ws.Range("K3").Value2 = 3.0
$ws.Range("M3").Value2 = 20.25989766666666
$ws.Range("N3").Value2 = 60.77969299999999
$ws.Range("O3").Value2 = 0.1424583290084953
$ws.Range("P3").Value2 = 0.1710942741446817
$ws.Range("Q3").Value2 = 456.9891806543353
$ws.Range("R3").Value2 = 2741.935083926012
$ws.Range("S3").Value2 = 0.009510931444827201
$ws.Range("T3").Value2 = 0.007982793508188807

$ws.Range("E4").Value2 = 2.0
$ws.Range("G4").Value2 = 22.556342
$ws.Range("H4").Value2 = 45.112684
$ws.Range("I4").Value2 = 0.0667629019027735
$ws.Range("J4").Value2 = 0.04665728030990886
$ws.Range("K4").Value2 = 3.0
$ws.Range("M4").Value2 = 22.68048566666667
$ws.Range("N4").Value2 = 68.041457
$ws.Range("O4").Value2 = 0.1594787961091443
$ws.Range("P4").Value2 = 0.1915360727004918
$ws.Range("Q4").Value2 = 511.5887914234314
$ws.Range("R4").Value2 = 3069.532748540588
$ws.Range("S4").Value2 = 0.01064726722020722
$ws.Range("T4").Value2 = 0.008936552233445928

$ws.Range("E5").Value2 = 2.0
$ws.Range("G5").Value2 = 22.556342
$ws.Range("H5").Value2 = 45.112684
$ws.Range("I5").Value2 = 0.0667629019027735
$ws.Range("J5").Value2 = 0.04665728030990886
$ws.Range("K5").Value2 = 3.0
$ws.Range("M5").Value2 = 10.83820733333333
$ws.Range("N5").Value2 = 32.514622
$ws.Range("O5").Value2 = 0.07620931416127522
$ws.Range("P5").Value2 = 0.09152836047030874
$ws.Range("Q5").Value2 = 244.4703112775747
$ws.Range("R5").Value2 = 1466.821867665448
$ws.Range("S5").Value2 = 0.005087954965426864
$ws.Range("T5").Value2 = 0.004270464370769576

$ws.Range("E6").Value2 = 2.0
$ws.Range("G6").Value2 = 22.556342
$ws.Range("H6").Value2 = 45.112684
$ws.Range("I6").Value2 = 0.0667629019027735
$ws.Range("J6").Value2 = 0.04665728030990886
$ws.Range("K6").Value2 = 3.0
$ws.Range("M6").Value2 = 17.02974533333333
$ws.Range("N6").Value2 = 51.089236
$ws.Range("O6").Value2 = 0.1197453759906399
$ws.Range("P6").Value2 = 0.1438157272368313
$ws.Range("Q6").Value2 = 384.1287599115707
$ws.Range("R6").Value2 = 2304.772559469424
$ws.Range("S6").Value2 = 0.007994548790573819
$ws.Range("T6").Value2 = 0.006710050698662231

$ws.Range("E7").Value2 = 2.0
$ws.Range("G7").Value2 = 22.556342
$ws.Range("H7").Value2 = 45.112684
$ws.Range("I7").Value2 = 0.0667629019027735
$ws.Range("J7").Value2 = 0.04665728030990886
$ws.Range("K7").Value2 = 2.0
$ws.Range("M7").Value2 = 41.965296
$ws.Range("N7").Value2 = 83.930592
$ws.Range("O7").Value2 = 0.2950807572114699
$ws.Range("P7").Value2 = 0.2362638408978708
$ws.Range("Q7").Value2 = 946.5835687072321
$ws.Range("R7").Value2 = 3786.334274828928
$ws.Range("S7").Value2 = 0.01970044764710549
$ws.Range("T7").Value2 = 0.01102342825186767

$ws.Range("E8").Value2 = 3.0
$ws.Range("G8").Value2 = 136.2606836666667
$ws.Range("H8").Value2 = 408.782051
$ws.Range("I8").Value2 = 0.4033082428366495
$ws.Range("J8").Value2 = 0.4227781867105593
$ws.Range("K8").Value2 = 2.0
$ws.Range("M8").Value2 = 29.4426765
$ws.Range("N8").Value2 = 58.88535299999999
$ws.Range("O8").Value2 = 0.2070274275189754
$ws.Range("P8").Value2 = 0.1657617245498156
$ws.Range("Q8").Value2 = 4011.8792288665
$ws.Range("R8").Value2 = 24071.275373199
$ws.Range("S8").Value2 = 0.0834958680116698
$ws.Range("T8").Value2 = 0.07008044133118624

$ws.Range("E9").Value2 = 3.0
$ws.Range("G9").Value2 = 136.2606836666667
$ws.Range("H9").Value2 = 408.782051
$ws.Range("I9").Value2 = 0.4033082428366495
$ws.Range("J9").Value2 = 0.4227781867105593
$ws.Range("K9").Value2 = 3.0
$ws.Range("M9").Value2 = 20.25989766666666
$ws.Range("N9").Value2 = 60.77969299999999
$ws.Range("O9").Value2 = 0.1424583290084953
$ws.Range("P9").Value2 = 0.1710942741446817
$ws.Range("Q9").Value2 = 2760.627507076704
$ws.Range("R9").Value2 = 24845.64756369034
$ws.Range("S9").Value2 = 0.05745461834986153
$ws.Range("T9").Value2 = 0.07233492697944786

$ws.Range("E10").Value2 = 3.0
$ws.Range("G10").Value2 = 136.2606836666667
$ws.Range("H10").Value2 = 408.782051
$ws.Range("I10").Value2 = 0.4033082428366495
$ws.Range("J10").Value2 = 0.4227781867105593
$ws.Range("K10").Value2 = 3.0
$ws.Range("M10").Value2 = 22.68048566666667
$ws.Range("N10").Value2 = 68.041457
$ws.Range("O10").Value2 = 0.1594787961091443
$ws.Range("P10").Value2 = 0.1915360727004918
$ws.Range("Q10").Value2 = 3090.458482832034
$ws.Range("R10").Value2 = 27814.1263454883
$ws.Range("S10").Value2 = 0.06431911302848328
$ws.Range("T10").Value2 = 0.08097727350597578

$ws.Range("E11").Value2 = 3.0
$ws.Range("G11").Value2 = 136.2606836666667
$ws.Range("H11").Value2 = 408.782051
$ws.Range("I11").Value2 = 0.4033082428366495
$ws.Range("J11").Value2 = 0.4227781867105593
$ws.Range("K11").Value2 = 3.0
$ws.Range("M11").Value2 = 10.83820733333333
$ws.Range("N11").Value2 = 32.514622
$ws.Range("O11").Value2 = 0.07620931416127522
$ws.Range("P11").Value2 = 0.09152836047030874
$ws.Range("Q11").Value2 = 1476.82154096108
$ws.Range("R11").Value2 = 13291.39386864972
$ws.Range("S11").Value2 = 0.0307358445821701
$ws.Range("T11").Value2 = 0.03869619427222756

$ws.Range("E12").Value2 = 3.0
$ws.Range("G12").Value2 = 136.2606836666667
$ws.Range("H12").Value2 = 408.782051
$ws.Range("I12").Value2 = 0.4033082428366495
$ws.Range("J12").Value2 = 0.4227781867105593
$ws.Range("K12").Value2 = 3.0
$ws.Range("M12").Value2 = 17.02974533333333
$ws.Range("N12").Value2 = 51.089236
$ws.Range("O12").Value2 = 0.1197453759906399
$ws.Range("P12").Value2 = 0.1438157272368313
$ws.Range("Q12").Value2 = 2320.484741789227
$ws.Range("R12").Value2 = 20884.36267610304
$ws.Range("S12").Value2 = 0.04829429717859889
$ws.Range("T12").Value2 = 0.06080215238164792

$ws.Range("E13").Value2 = 3.0
$ws.Range("G13").Value2 = 136.2606836666667
$ws.Range("H13").Value2 = 408.782051
$ws.Range("I13").Value2 = 0.4033082428366495
$ws.Range("J13").Value2 = 0.4227781867105593
$ws.Range("K13").Value2 = 2.0
$ws.Range("M13").Value2 = 41.965296
$ws.Range("N13").Value2 = 83.930592
$ws.Range("O13").Value2 = 0.2950807572114699
$ws.Range("P13").Value2 = 0.2362638408978708
$ws.Range("Q13").Value2 = 5718.219923234033
$ws.Range("R13").Value2 = 34309.3195394042
$ws.Range("S13").Value2 = 0.1190085016858659
$ws.Range("T13").Value2 = 0.0998871982400739

$ws.Range("E14").Value2 = 3.0
$ws.Range("G14").Value2 = 60.036368
$ws.Range("H14").Value2 = 180.109104
$ws.Range("I14").Value2 = 0.1776973477074789
$ws.Range("J14").Value2 = 0.1862757922293989
$ws.Range("K14").Value2 = 2.0
$ws.Range("M14").Value2 = 29.4426765
$ws.Range("N14").Value2 = 58.88535299999999
$ws.Range("O14").Value2 = 0.2070274275189754
$ws.Range("P14").Value2 = 0.1657617245498156
$ws.Range("Q14").Value2 = 1767.631361258952
$ws.Range("R14").Value2 = 10605.78816755371
$ws.Range("S14").Value2 = 0.03678822477282426
$ws.Range("T14").Value2 = 0.0308773965618283

$ws.Range("E15").Value2 = 3.0
$ws.Range("G15").Value2 = 60.036368
$ws.Range("H15").Value2 = 180.109104
$ws.Range("I15").Value2 = 0.1776973477074789
$ws.Range("J15").Value2 = 0.1862757922293989
$ws.Range("K15").Value2 = 3.0
$ws.Range("M15").Value2 = 20.25989766666666
$ws.Range("N15").Value2 = 60.77969299999999
$ws.Range("O15").Value2 = 0.1424583290084953
$ws.Range("P15").Value2 = 0.1710942741446817
$ws.Range("Q15").Value2 = 1216.330671958341
$ws.Range("R15").Value2 = 10946.97604762507
$ws.Range("S15").Value2 = 0.02531446722364901
$ws.Range("T15").Value2 = 0.03187072146221455

$ws.Range("E16").Value2 = 3.0
$ws.Range("G16").Value2 = 60.036368
$ws.Range("H16").Value2 = 180.109104
$ws.Range("I16").Value2 = 0.1776973477074789
$ws.Range("J16").Value2 = 0.1862757922293989
$ws.Range("K16").Value2 = 3.0
$ws.Range("M16").Value2 = 22.68048566666667
$ws.Range("N16").Value2 = 68.041457
$ws.Range("O16").Value2 = 0.1594787961091443
$ws.Range("P16").Value2 = 0.1915360727004918
$ws.Range("Q16").Value2 = 1361.653983902725
$ws.Range("R16").Value2 = 12254.88585512453
$ws.Range("S16").Value2 = 0.02833895908417675
$ws.Range("T16").Value2 = 0.03567853368279185

$ws.Range("E17").Value2 = 3.0
$ws.Range("G17").Value2 = 60.036368
$ws.Range("H17").Value2 = 180.109104
$ws.Range("I17").Value2 = 0.1776973477074789
$ws.Range("J17").Value2 = 0.1862757922293989
$ws.Range("K17").Value2 = 3.0
$ws.Range("M17").Value2 = 10.83820733333333
$ws.Range("N17").Value2 = 32.514622
$ws.Range("O17").Value2 = 0.07620931416127522
$ws.Range("P17").Value2 = 0.09152836047030874
$ws.Range("Q17").Value2 = 650.6866039242988
$ws.Range("R17").Value2 = 5856.179435318689
$ws.Range("S17").Value2 = 0.01354219299706462
$ws.Range("T17").Value2 = 0.01704951785806476

$ws.Range("E18").Value2 = 3.0
$ws.Range("G18").Value2 = 60.036368
$ws.Range("H18").Value2 = 180.109104
$ws.Range("I18").Value2 = 0.1776973477074789
$ws.Range("J18").Value2 = 0.1862757922293989
$ws.Range("K18").Value2 = 3.0
$ws.Range("M18").Value2 = 17.02974533333333
$ws.Range("N18").Value2 = 51.089236
$ws.Range("O18").Value2 = 0.1197453759906399
$ws.Range("P18").Value2 = 0.1438157272368313
$ws.Range("Q18").Value2 = 1022.404057778283
$ws.Range("R18").Value2 = 9201.636520004544
$ws.Range("S18").Value2 = 0.02127843571377153
$ws.Range("T18").Value2 = 0.02678938852608789

$ws.Range("E19").Value2 = 3.0
$ws.Range("G19").Value2 = 60.036368
$ws.Range("H19").Value2 = 180.109104
$ws.Range("I19").Value2 = 0.1776973477074789
$ws.Range("J19").Value2 = 0.1862757922293989
$ws.Range("K19").Value2 = 2.0
$ws.Range("M19").Value2 = 41.965296
$ws.Range("N19").Value2 = 83.930592
$ws.Range("O19").Value2 = 0.2950807572114699
$ws.Range("P19").Value2 = 0.2362638408978708
$ws.Range("Q19").Value2 = 2519.443953884928
$ws.Range("R19").Value2 = 15116.66372330957
$ws.Range("S19").Value2 = 0.05243506791599273
$ws.Range("T19").Value2 = 0.04401023413841154

$ws.Range("E20").Value2 = 3.0
$ws.Range("G20").Value2 = 67.73010366666666
$ws.Range("H20").Value2 = 203.190311
$ws.Range("I20").Value2 = 0.200469485121406
$ws.Range("J20").Value2 = 0.2101472680407257
$ws.Range("K20").Value2 = 2.0
$ws.Range("M20").Value2 = 29.4426765
$ws.Range("N20").Value2 = 58.88535299999999
$ws.Range("O20").Value2 = 0.2070274275189754
$ws.Range("P20").Value2 = 0.1657617245498156
$ws.Range("Q20").Value2 = 1994.15553156913
$ws.Range("R20").Value2 = 11964.93318941478
$ws.Range("S20").Value2 = 0.0415026818007382
$ws.Range("T20").Value2 = 0.03483437355986304

$ws.Range("E21").Value2 = 3.0
$ws.Range("G21").Value2 = 67.73010366666666
$ws.Range("H21").Value2 = 203.190311
$ws.Range("I21").Value2 = 0.200469485121406
$ws.Range("J21").Value2 = 0.2101472680407257
$ws.Range("K21").Value2 = 3.0
$ws.Range("M21").Value2 = 20.25989766666666
$ws.Range("N21").Value2 = 60.77969299999999
$ws.Range("O21").Value2 = 0.1424583290084953
$ws.Range("P21").Value2 = 0.1710942741446817
$ws.Range("Q21").Value2 = 1372.204969239391
$ws.Range("R21").Value2 = 12349.84472315452
$ws.Range("S21").Value2 = 0.02855854786758891
$ws.Range("T21").Value2 = 0.03595499428891583

$ws.Range("E22").Value2 = 3.0
$ws.Range("G22").Value2 = 67.73010366666666
$ws.Range("H22").Value2 = 203.190311
$ws.Range("I22").Value2 = 0.200469485121406
$ws.Range("J22").Value2 = 0.2101472680407257
$ws.Range("K22").Value2 = 3.0
$ws.Range("M22").Value2 = 22.68048566666667
$ws.Range("N22").Value2 = 68.041457
$ws.Range("O22").Value2 = 0.1594787961091443
$ws.Range("P22").Value2 = 0.1915360727004918
$ws.Range("Q22").Value2 = 1536.151645413681
$ws.Range("R22").Value2 = 13825.36480872313
$ws.Range("S22").Value2 = 0.03197063214378185
$ws.Range("T22").Value2 = 0.04025078240925817

$ws.Range("E23").Value2 = 3.0
$ws.Range("G23").Value2 = 67.73010366666666
$ws.Range("H23").Value2 = 203.190311
$ws.Range("I23").Value2 = 0.200469485121406
$ws.Range("J23").Value2 = 0.2101472680407257
$ws.Range("K23").Value2 = 3.0
$ws.Range("M23").Value2 = 10.83820733333333
$ws.Range("N23").Value2 = 32.514622
$ws.Range("O23").Value2 = 0.07620931416127522
$ws.Range("P23").Value2 = 0.09152836047030874
$ws.Range("Q23").Value2 = 734.0729062474936
$ws.Range("R23").Value2 = 6606.656156227443
$ws.Range("S23").Value2 = 0.01527764197136632
$ws.Range("T23").Value2 = 0.01923443490108213

$ws.Range("E24").Value2 = 3.0
$ws.Range("G24").Value2 = 67.73010366666666
$ws.Range("H24").Value2 = 203.190311
$ws.Range("I24").Value2 = 0.200469485121406
$ws.Range("J24").Value2 = 0.2101472680407257
$ws.Range("K24").Value2 = 3.0
$ws.Range("M24").Value2 = 17.02974533333333
$ws.Range("N24").Value2 = 51.089236
$ws.Range("O24").Value2 = 0.1197453759906399
$ws.Range("P24").Value2 = 0.1438157272368313
$ws.Range("Q24").Value2 = 1153.4264168436
$ws.Range("R24").Value2 = 10380.8377515924
$ws.Range("S24").Value2 = 0.02400529387051275
$ws.Range("T24").Value2 = 0.03022248218011028

$ws.Range("E25").Value2 = 3.0
$ws.Range("G25").Value2 = 67.73010366666666
$ws.Range("H25").Value2 = 203.190311
$ws.Range("I25").Value2 = 0.200469485121406
$ws.Range("J25").Value2 = 0.2101472680407257
$ws.Range("K25").Value2 = 2.0
$ws.Range("M25").Value2 = 41.965296
$ws.Range("N25").Value2 = 83.930592
$ws.Range("O25").Value2 = 0.2950807572114699
$ws.Range("P25").Value2 = 0.2362638408978708
$ws.Range("Q25").Value2 = 2842.313848482352
$ws.Range("R25").Value2 = 17053.88309089411
$ws.Range("S25").Value2 = 0.05915468746741798
$ws.Range("T25").Value2 = 0.04965020070149623

$ws.Range("E26").Value2 = 3.0
$ws.Range("G26").Value2 = 27.15284733333333
$ws.Range("H26").Value2 = 81.458542
$ws.Range("I26").Value2 = 0.0803677689802858
$ws.Range("J26").Value2 = 0.08424757054425056
$ws.Range("K26").Value2 = 2.0
$ws.Range("M26").Value2 = 29.4426765
$ws.Range("N26").Value2 = 58.88535299999999
$ws.Range("O26").Value2 = 0.2070274275189754
$ws.Range("P26").Value2 = 0.1657617245498156
$ws.Range("Q26").Value2 = 799.4525000892208
$ws.Range("R26").Value2 = 4796.715000535325
$ws.Range("S26").Value2 = 0.01663833246742788
$ws.Range("T26").Value2 = 0.01396502258254722

$ws.Range("E27").Value2 = 3.0
$ws.Range("G27").Value2 = 27.15284733333333
$ws.Range("H27").Value2 = 81.458542
$ws.Range("I27").Value2 = 0.0803677689802858
$ws.Range("J27").Value2 = 0.08424757054425056
$ws.Range("K27").Value2 = 3.0
$ws.Range("M27").Value2 = 20.25989766666666
$ws.Range("N27").Value2 = 60.77969299999999
$ws.Range("O27").Value2 = 0.1424583290084953
$ws.Range("P27").Value2 = 0.1710942741446817
$ws.Range("Q27").Value2 = 550.113908331956
$ws.Range("R27").Value2 = 4951.025174987605
$ws.Range("S27").Value2 = 0.0114490580750723
$ws.Range("T27").Value2 = 0.01441427693072142

$ws.Range("E28").Value2 = 3.0
$ws.Range("G28").Value2 = 27.15284733333333
$ws.Range("H28").Value2 = 81.458542
$ws.Range("I28").Value2 = 0.0803677689802858
$ws.Range("J28").Value2 = 0.08424757054425056
$ws.Range("K28").Value2 = 3.0
$ws.Range("M28").Value2 = 22.68048566666667
$ws.Range("N28").Value2 = 68.041457
$ws.Range("O28").Value2 = 0.1594787961091443
$ws.Range("P28").Value2 = 0.1915360727004918
$ws.Range("Q28").Value2 = 615.8397647528548
$ws.Range("R28").Value2 = 5542.557882775694
$ws.Range("S28").Value2 = 0.01281695504295381
$ws.Range("T28").Value2 = 0.01613644879660339

$ws.Range("E29").Value2 = 3.0
$ws.Range("G29").Value2 = 27.15284733333333
$ws.Range("H29").Value2 = 81.458542
$ws.Range("I29").Value2 = 0.0803677689802858
$ws.Range("J29").Value2 = 0.08424757054425056
$ws.Range("K29").Value2 = 3.0
$ws.Range("M29").Value2 = 10.83820733333333
$ws.Range("N29").Value2 = 32.514622
$ws.Range("O29").Value2 = 0.07620931416127522
$ws.Range("P29").Value2 = 0.09152836047030874
$ws.Range("Q29").Value2 = 294.2881890890138
$ws.Range("R29").Value2 = 2648.593701801124
$ws.Range("S29").Value2 = 0.00612477255465939
$ws.Range("T29").Value2 = 0.00771104200552193

$ws.Range("E30").Value2 = 3.0
$ws.Range("G30").Value2 = 27.15284733333333
$ws.Range("H30").Value2 = 81.458542
$ws.Range("I30").Value2 = 0.0803677689802858
$ws.Range("J30").Value2 = 0.08424757054425056
$ws.Range("K30").Value2 = 3.0
$ws.Range("M30").Value2 = 17.02974533333333
$ws.Range("N30").Value2 = 51.089236
$ws.Range("O30").Value2 = 0.1197453759906399
$ws.Range("P30").Value2 = 0.1438157272368313
$ws.Range("Q30").Value2 = 462.4060751615457
$ws.Range("R30").Value2 = 4161.654676453912
$ws.Range("S30").Value2 = 0.009623668714073209
$ws.Range("T30").Value2 = 0.01211612562575764

$ws.Range("E31").Value2 = 3.0
$ws.Range("G31").Value2 = 27.15284733333333
$ws.Range("H31").Value2 = 81.458542
$ws.Range("I31").Value2 = 0.0803677689802858
$ws.Range("J31").Value2 = 0.08424757054425056
$ws.Range("K31").Value2 = 2.0
$ws.Range("M31").Value2 = 41.965296
$ws.Range("N31").Value2 = 83.930592
$ws.Range("O31").Value2 = 0.2950807572114699
$ws.Range("P31").Value2 = 0.2362638408978708
$ws.Range("Q31").Value2 = 1139.477275586144
$ws.Range("R31").Value2 = 6836.863653516863
$ws.Range("S31").Value2 = 0.02371498212609922
$ws.Range("T31").Value2 = 0.01990465460309896

$ws.Range("E32").Value2 = 2.0
$ws.Range("G32").Value2 = 24.1210785
$ws.Range("H32").Value2 = 48.242157
$ws.Range("I32").Value2 = 0.0713942534514062
$ws.Range("J32").Value2 = 0.04989390216515674
$ws.Range("K32").Value2 = 2.0
$ws.Range("M32").Value2 = 29.4426765
$ws.Range("N32").Value2 = 58.88535299999999
$ws.Range("O32").Value2 = 0.2070274275189754
$ws.Range("P32").Value2 = 0.1657617245498156
$ws.Range("Q32").Value2 = 710.1891111066052
$ws.Range("R32").Value2 = 2840.756444426421
$ws.Range("S32").Value2 = 0.01478056863168236
$ws.Range("T32").Value2 = 0.00827049926741616

$ws.Range("E33").Value2 = 2.0
$ws.Range("G33").Value2 = 24.1210785
$ws.Range("H33").Value2 = 48.242157
$ws.Range("I33").Value2 = 0.0713942534514062
$ws.Range("J33").Value2 = 0.04989390216515674
$ws.Range("K33").Value2 = 3.0
$ws.Range("M33").Value2 = 20.25989766666666
$ws.Range("N33").Value2 = 60.77969299999999
$ws.Range("O33").Value2 = 0.1424583290084953
$ws.Range("P33").Value2 = 0.1710942741446817
$ws.Range("Q33").Value2 = 488.6905820196334
$ws.Range("R33").Value2 = 2932.143492117801
$ws.Range("S33").Value2 = 0.01017070604749632
$ws.Range("T33").Value2 = 0.008536560975193256

$ws.Range("E34").Value2 = 2.0
$ws.Range("G34").Value2 = 24.1210785
$ws.Range("H34").Value2 = 48.242157
$ws.Range("I34").Value2 = 0.0713942534514062
$ws.Range("J34").Value2 = 0.04989390216515674
$ws.Range("K34").Value2 = 3.0
$ws.Range("M34").Value2 = 22.68048566666667
$ws.Range("N34").Value2 = 68.041457
$ws.Range("O34").Value2 = 0.1594787961091443
$ws.Range("P34").Value2 = 0.1915360727004918
$ws.Range("Q34").Value2 = 547.0777751837915
$ws.Range("R34").Value2 = 3282.466651102749
$ws.Range("S34").Value2 = 0.01138586958954138
$ws.Range("T34").Value2 = 0.009556482072416687

$ws.Range("E35").Value2 = 2.0
$ws.Range("G35").Value2 = 24.1210785
$ws.Range("H35").Value2 = 48.242157
$ws.Range("I35").Value2 = 0.0713942534514062
$ws.Range("J35").Value2 = 0.04989390216515674
$ws.Range("K35").Value2 = 3.0
$ws.Range("M35").Value2 = 10.83820733333333
$ws.Range("N35").Value2 = 32.514622
$ws.Range("O35").Value2 = 0.07620931416127522
$ws.Range("P35").Value2 = 0.09152836047030874
$ws.Range("Q35").Value2 = 261.4292498866091
$ws.Range("R35").Value2 = 1568.575499319654
$ws.Range("S35").Value2 = 0.005440907090587922
$ws.Range("T35").Value2 = 0.004566707062642784

$ws.Range("E36").Value2 = 2.0
$ws.Range("G36").Value2 = 24.1210785
$ws.Range("H36").Value2 = 48.242157
$ws.Range("I36").Value2 = 0.0713942534514062
$ws.Range("J36").Value2 = 0.04989390216515674
$ws.Range("K36").Value2 = 3.0
$ws.Range("M36").Value2 = 17.02974533333333
$ws.Range("N36").Value2 = 51.089236
$ws.Range("O36").Value2 = 0.1197453759906399
$ws.Range("P36").Value2 = 0.1438157272368313
$ws.Range("Q36").Value2 = 410.775824020342
$ws.Range("R36").Value2 = 2464.654944122052
$ws.Range("S36").Value2 = 0.008549131723109674
$ws.Range("T36").Value2 = 0.007175527824565328

$ws.Range("E37").Value2 = 2.0
$ws.Range("G37").Value2 = 24.1210785
$ws.Range("H37").Value2 = 48.242157
$ws.Range("I37").Value2 = 0.0713942534514062
$ws.Range("J37").Value2 = 0.04989390216515674
$ws.Range("K37").Value2 = 2.0
$ws.Range("M37").Value2 = 41.965296
$ws.Range("N37").Value2 = 83.930592
$ws.Range("O37").Value2 = 0.2950807572114699
$ws.Range("P37").Value2 = 0.2362638408978708
$ws.Range("Q37").Value2 = 1012.248199091736
$ws.Range("R37").Value2 = 4048.992796366944
$ws.Range("S37").Value2 = 0.02106707036898854
$ws.Range("T37").Value2 = 0.01178812496292253

Write-Output "Edit complete"